# Auto commit at 2025-12-27  8:07:35.66
# Updates the raw metric figures on the "Metrics" sheet. Dependent formulas
# on the "today" sheet (and the TODAY()-1 cell) recalc automatically.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value = 358718.41000000003
$metrics.Range("B3").Value = 307439.72999999992
$metrics.Range("B4").Value = 109722.98999999999
$metrics.Range("B5").Value = 14566
$metrics.Range("B6").Value = 5561425.5200000005
$metrics.Range("B7").Value = 4707792.6900000004
$metrics.Range("B8").Value = 1641679.87
$metrics.Range("B9").Value = 217273
$metrics.Range("B10").Value = 34026806.509999998
$metrics.Range("B11").Value = 31983067.850000001
$metrics.Range("B12").Value = 11923401.909999995
$metrics.Range("B13").Value = 1314903

# Update the recorded cursor/selection on each sheet that moved.
$metrics.Range("E15").Select()

$today = $wb.Worksheets.Item("today")
$today.Range("F7").Select()
